{"js": "// Localize a handful of strings in the (Chinese) resume document:\n//  - switch a \" - \" separator for a full-width colon in the title\n//  - normalize full-width parentheses to ASCII ones in the phone number\n//  - translate company names to English and update employment dates/wording\n//  - translate the publisher name in the bibliography entry\n//\n// Each change is scoped to a single paragraph, so we find-and-replace the\n// exact (unique) original text of each paragraph with its edited text.\n\nconst replacements = [\n  [\"\u7b80\u5386 - Patti Fernandez\", \"\u7b80\u5386\uff1aPatti Fernandez\"],\n  [\"\u7535\u8bdd\uff1a\uff08123\uff09 456-7890\", \"\u7535\u8bdd\uff1a(123) 456-7890\"],\n  [\"ABC \u5de5\u4f5c\u5ba4\uff1a\u9996\u5e2d\u52a8\u753b\u5e08\uff082018 \u5e74 1 \u6708 - \u6f14\u793a\uff09\", \"ABC Studios\uff1a\u9996\u5e2d\u52a8\u753b\u5e08\uff082018 \u5e74 1 \u6708 - \u4eca\uff09\"],\n  [\"XYZ \u5a92\u4f53\uff1a\u9ad8\u7ea7\u52a8\u753b\u5e08 \uff082015 \u5e74 6 \u6708 - 2017 \u5e74 12 \u6708\uff09\", \"XYZ Media\uff1a\u9ad8\u7ea7\u52a8\u753b\u5e08\uff082015 \u5e74 6 \u6708 - 2017 \u5e74 12 \u6708\uff09\"],\n  [\"MNO \u5a31\u4e50\uff1a \u521d\u7ea7\u52a8\u753b\u5e08 \uff082012 \u5e74 9 \u6708 - 2015 \u5e74 5 \u6708\uff09\", \"MNO Entertainment\uff1a\u521d\u7ea7\u52a8\u753b\u5e08 \uff082012 \u5e74 9 \u6708 - 2015 \u5e74 5 \u6708\uff09\"],\n  [\"\u7ebd\u7ea6\uff1a\u4f01\u9e45\u4e66\u3002\", \"\u7ebd\u7ea6\uff1aPenguin Books \u51fa\u7248\u793e\u3002\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Localize a handful of strings in the (Chinese) resume document:\n#  - switch a \" - \" separator for a full-width colon in the title\n#  - normalize full-width parentheses to ASCII ones in the phone number\n#  - translate company names to English and update employment dates/wording\n#  - translate the publisher name in the bibliography entry\n#\n# Each change is a unique, whole-document Find & Replace (wdReplaceAll),\n# scoped by the exact original paragraph text so each edit lands in only\n# the one spot the diff targets.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"\u7b80\u5386 - Patti Fernandez\", \"\u7b80\u5386\uff1aPatti Fernandez\"),\n    @(\"\u7535\u8bdd\uff1a\uff08123\uff09 456-7890\", \"\u7535\u8bdd\uff1a(123) 456-7890\"),\n    @(\"ABC \u5de5\u4f5c\u5ba4\uff1a\u9996\u5e2d\u52a8\u753b\u5e08\uff082018 \u5e74 1 \u6708 - \u6f14\u793a\uff09\", \"ABC Studios\uff1a\u9996\u5e2d\u52a8\u753b\u5e08\uff082018 \u5e74 1 \u6708 - \u4eca\uff09\"),\n    @(\"XYZ \u5a92\u4f53\uff1a\u9ad8\u7ea7\u52a8\u753b\u5e08 \uff082015 \u5e74 6 \u6708 - 2017 \u5e74 12 \u6708\uff09\", \"XYZ Media\uff1a\u9ad8\u7ea7\u52a8\u753b\u5e08\uff082015 \u5e74 6 \u6708 - 2017 \u5e74 12 \u6708\uff09\"),\n    @(\"MNO \u5a31\u4e50\uff1a \u521d\u7ea7\u52a8\u753b\u5e08 \uff082012 \u5e74 9 \u6708 - 2015 \u5e74 5 \u6708\uff09\", \"MNO Entertainment\uff1a\u521d\u7ea7\u52a8\u753b\u5e08 \uff082012 \u5e74 9 \u6708 - 2015 \u5e74 5 \u6708\uff09\"),\n    @(\"\u7ebd\u7ea6\uff1a\u4f01\u9e45\u4e66\u3002\", \"\u7ebd\u7ea6\uff1aPenguin Books \u51fa\u7248\u793e\u3002\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $oldText\n    $range.Find.Replacement.Text = $newText\n    $range.Find.Forward = $true\n    $range.Find.MatchCase = $true\n    $range.Find.MatchWholeWord = $false\n    $range.Find.MatchWildcards = $false\n    $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
